$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = 70939; B = "Maysa Lima";             C = "Atendimento ao Cliente"; D = "Outros";              E = 7; F = 45084; G = 4379.68 },
    @{ Row = 3;  A = 14783; B = "Fernando Almeida";        C = "P&D";                    D = "Problemas pessoais";  E = 3; F = 45095; G = 3800.78 },
    @{ Row = 4;  A = 89144; B = "Sra. Ana Beatriz Dias";   C = "Marketing";               D = "Consulta médica";     E = 1; F = 45095; G = 11627.34 },
    @{ Row = 5;  A = 63925; B = "Sofia Costa";             C = "Financeiro";              D = "Doença";              E = 8; F = 45098; G = 8847.219999999999 },
    @{ Row = 6;  A = 20899; B = "Isadora Campos";          C = "TI";                      D = "Consulta médica";     E = 1; F = 45086; G = 7820.2 },
    @{ Row = 7;  A = 88469; B = "Maria Fernanda da Paz";   C = "Engenharia";              D = "Problemas pessoais";  E = 6; F = 45100; G = 5302.59 },
    @{ Row = 8;  A = 79083; B = "Dr. João Pedro Martins";  C = "P&D";                    D = "Outros";              E = 8; F = 45097; G = 3889.21 },
    @{ Row = 9;  A = 46864; B = "Ana Clara Rocha";         C = "Vendas";                  D = "Problemas pessoais";  E = 4; F = 45079; G = 7733.77 },
    @{ Row = 10; A = 49648; B = "Alice Lopes";             C = "Vendas";                  D = "Viagem de negócios";  E = 8; F = 45085; G = 4690.48 },
    @{ Row = 11; A = 97219; B = "Yasmin Correia";          C = "Atendimento ao Cliente"; D = "Problemas pessoais";  E = 6; F = 45101; G = 3915.61 }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value = $rec.A
    $ws.Cells.Item($r, 2).Value = $rec.B
    $ws.Cells.Item($r, 3).Value = $rec.C
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $rec.E
    $ws.Cells.Item($r, 6).Value = $rec.F
    $ws.Cells.Item($r, 7).Value = $rec.G
}
